# Spring 24 week 8 inputs: append 33 new matchup rows to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(5, 0, 5, 2),
    @(4, 3, 3, 0),
    @(4, 0, 5, 3),
    @(3, 0, 3, 3),
    @(6, 0, 7, 2),
    @(2, 0, 3, 3),
    @(6, 2, 5, 0),
    @(5, 0, 7, 3),
    @(3, 2, 4, 1),
    @(6, 2, 4, 1),
    @(4, 0, 4, 3),
    @(6, 2, 3, 1),
    @(3, 0, 5, 3),
    @(4, 2, 4, 0),
    @(5, 2, 7, 1),
    @(4, 0, 4, 2),
    @(6, 2, 6, 1),
    @(5, 2, 4, 0),
    @(4, 0, 3, 2),
    @(4, 2, 7, 0),
    @(3, 3, 3, 0),
    @(4, 0, 5, 2),
    @(3, 1, 3, 2),
    @(4, 1, 4, 2),
    @(3, 0, 2, 2),
    @(3, 0, 4, 3),
    @(7, 2, 4, 1),
    @(6, 1, 6, 2),
    @(4, 1, 5, 2),
    @(4, 3, 4, 0),
    @(3, 0, 4, 3),
    @(6, 2, 5, 1),
    @(5, 1, 4, 2)
)

$startRow = 2672
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newData[$i]
    for ($col = 1; $col -le 4; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}

# Move selection to the next empty row, matching the saved workbook state
$ws.Range("A2705").Select()
